# itemList.xlsx - refresh the sample wishlist/cart item text (new screenshot +
# allure-report fixture data) and leave the cart tab focused, matching the
# state the workbook was left in after the author's manual Excel edit.

$wb = $excel.ActiveWorkbook

$wsWish = $wb.Worksheets.Item("wishlistItem")
$wsCart = $wb.Worksheets.Item("cartItem")

# --- wishlistItem sheet: replace the three wish-listed product names ---
$wsWish.Range("A1").Value = "The Great British Barbeque Company Instand BBQ With Stand"
$wsWish.Range("A2").Value = "Doff Grow Bag Multipurpose Potting Soil"
$wsWish.Range("A3").Value = "Hanging Basket with Liner"

# --- cartItem sheet: replace the three cart product names ---
$wsCart.Range("A1").Value = "Solar Powered LED Outdoor Stake Light - Orange"
$wsCart.Range("A2").Value = "Wilson & Gregory Heavy Duty Garden Bracket"
$wsCart.Range("A3").Value = "Pepco Solar Powered Stake Lights (Pack of 12)"

# The new cartItem text wraps onto a second line, so rows 2 and 3 grow to
# match row 1's wrapped height.
$wsCart.Rows.Item(2).RowHeight = 28.8
$wsCart.Rows.Item(3).RowHeight = 28.8

# Restore each sheet's own last-used cell selection.
[void]$wsWish.Range("A9").Select()
[void]$wsCart.Range("A3").Select()

# cartItem is the tab left in front (selecting on it last also makes it the
# active sheet/tab, matching tabSelected + workbook activeTab moving to it).
[void]$wsCart.Activate()
